$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 156; shifts rows 156-196 down to 157-197.
$ws.Rows.Item(156).Insert()

# Populate the new row 156 with the new weekly record.
# Non-date / non-changed columns mirror the existing data for this
# market/category/variety/quality combination.
$ws.Cells.Item(156, 1).Value = 7
$ws.Cells.Item(156, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(156, 3).Value = "Ñuble"
$ws.Cells.Item(156, 4).Value = 44642
$ws.Cells.Item(156, 5).Value = 16
$ws.Cells.Item(156, 6).Value = 100112006
$ws.Cells.Item(156, 7).Value = "Repollo"
$ws.Cells.Item(156, 8).Value = "Crespo record"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 150
$ws.Cells.Item(156, 11).Value = 1300
$ws.Cells.Item(156, 12).Value = 1300
$ws.Cells.Item(156, 13).Value = 1300
$ws.Cells.Item(156, 14).Value = '$/unidad'
$ws.Cells.Item(156, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(156, 16).Value = 1300
$ws.Cells.Item(156, 17).Value = 1
$ws.Cells.Item(156, 18).Value = "Hortaliza"
